$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-05 Wednesday" "2025-02-06 Thursday"

Replace-Text "124÷5=24, 4" "959÷7=137, 0"
Replace-Text "759÷8=94, 7" "764÷5=152, 4"
Replace-Text "787÷8=98, 3" "408÷4=102, 0"
Replace-Text "887÷7=126, 5" "363÷2=181, 1"
Replace-Text "101÷7=14, 3" "816÷2=408, 0"

Replace-Text "873÷6=145, 3" "573÷6=95, 3"
Replace-Text "886÷3=295, 1" "957÷9=106, 3"
Replace-Text "769÷2=384, 1" "697÷2=348, 1"
Replace-Text "924÷6=154, 0" "578÷2=289, 0"
Replace-Text "299÷7=42, 5" "984÷8=123, 0"

Replace-Text "407÷6=67, 5" "782÷8=97, 6"
Replace-Text "930÷9=103, 3" "893÷5=178, 3"
Replace-Text "293÷2=146, 1" "260÷7=37, 1"
Replace-Text "134÷6=22, 2" "330÷6=55, 0"
Replace-Text "445÷2=222, 1" "565÷4=141, 1"

Replace-Text "724÷9=80, 4" "796÷9=88, 4"
Replace-Text "772÷9=85, 7" "857÷7=122, 3"
Replace-Text "728÷8=91, 0" "265÷3=88, 1"
Replace-Text "910÷9=101, 1" "759÷3=253, 0"
Replace-Text "988÷7=141, 1" "694÷9=77, 1"

Replace-Text "708÷4=177, 0" "377÷8=47, 1"
Replace-Text "431÷4=107, 3" "797÷2=398, 1"
Replace-Text "432÷9=48, 0" "464÷7=66, 2"
Replace-Text "472÷2=236, 0" "347÷2=173, 1"
Replace-Text "403÷6=67, 1" "940÷7=134, 2"
